$d = $word.ActiveDocument
$d.Content.Find.Execute("fields no listed below", $true, $false, $false, $false, $false, $true, 1, $false, "fields not listed below", 2)
